$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "30.636.02"
$ws.Cells.Item(2,5).Value = "  +0.69%  "

# Row 3
$ws.Cells.Item(3,4).Value = "1.886.74"
$ws.Cells.Item(3,5).Value = "  +0.95%  "

# Row 4
$ws.Cells.Item(4,4).Value = "'1.003"
$ws.Cells.Item(4,5).Value = "  +0.25%  "

# Row 5
$ws.Cells.Item(5,4).Value = "'247.18"
$ws.Cells.Item(5,5).Value = "  +0.09%  "

# Row 6
$ws.Cells.Item(6,4).Value = "'1.002"
$ws.Cells.Item(6,5).Value = "  +0.18%  "

# Row 7
$ws.Cells.Item(7,4).Value = "'0.4716"
$ws.Cells.Item(7,5).Value = "  -0.41%  "

# Row 8
$ws.Cells.Item(8,4).Value = "'0.2918"
$ws.Cells.Item(8,5).Value = "  -0.04%  "

# Row 9
$ws.Cells.Item(9,4).Value = "'0.06515"
$ws.Cells.Item(9,5).Value = "  +0.45%  "

# Row 10
$ws.Cells.Item(10,4).Value = "'22.36"
$ws.Cells.Item(10,5).Value = "  +0.81%  "

# Row 11
$ws.Cells.Item(11,4).Value = "'0.07794"

# Row 12
$ws.Cells.Item(12,2).Value = "WrappedEther"
$ws.Cells.Item(12,3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12,4).Value = "1.890.49"
$ws.Cells.Item(12,5).Value = "  +1.06%  "

# Row 13
$ws.Cells.Item(13,2).Value = "Litecoin"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(13,4).Value = "'96.72"
$ws.Cells.Item(13,5).Value = "  -0.95%  "

# Row 14
$ws.Cells.Item(14,4).Value = "'0.7377"
$ws.Cells.Item(14,5).Value = "  -0.68%  "

# Row 15
$ws.Cells.Item(15,4).Value = "'5.230"
$ws.Cells.Item(15,5).Value = "  +1.68%  "

# Row 16
$ws.Cells.Item(16,4).Value = "'282.89"
$ws.Cells.Item(16,5).Value = "  +3.54%  "

# Row 17
$ws.Cells.Item(17,4).Value = "31.139.89"
$ws.Cells.Item(17,5).Value = "  +2.32%  "

# Row 18
$ws.Cells.Item(18,4).Value = "'13.22"
$ws.Cells.Item(18,5).Value = "  -1.51%  "

# Row 19
$ws.Cells.Item(19,2).Value = "Dai"
$ws.Cells.Item(19,3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(19,4).Value = "'1.003"
$ws.Cells.Item(19,5).Value = "  +0.31%  "

# Row 20
$ws.Cells.Item(20,2).Value = "ShibaInu"
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(20,4).Value = "'0.000007503"
$ws.Cells.Item(20,5).Value = "  -0.37%  "

# Row 21
$ws.Cells.Item(21,4).Value = "2.142.52"
$ws.Cells.Item(21,5).Value = "  +1.13%  "

# Row 22
$ws.Cells.Item(22,4).Value = "'5.308"
$ws.Cells.Item(22,5).Value = "  +1.37%  "

# Row 23
$ws.Cells.Item(23,4).Value = "'1.003"
$ws.Cells.Item(23,5).Value = "  +0.34%  "

# Row 24
$ws.Cells.Item(24,4).Value = "'6.249"
$ws.Cells.Item(24,5).Value = "  +1.21%  "

# Row 25
$ws.Cells.Item(25,4).Value = "'9.196"
$ws.Cells.Item(25,5).Value = "  -0.92%  "

# Row 26
$ws.Cells.Item(26,4).Value = "'164.55"
$ws.Cells.Item(26,5).Value = "  +0.73%  "

# Row 27
$ws.Cells.Item(27,4).Value = "'18.91"
$ws.Cells.Item(27,5).Value = "  +0.71%  "

# Row 28
$ws.Cells.Item(28,4).Value = "'1.910"
$ws.Cells.Item(28,5).Value = "  -0.81%  "

# Row 29
$ws.Cells.Item(29,4).Value = "'1.365"
$ws.Cells.Item(29,5).Value = "  +0.03%  "

# Row 30
$ws.Cells.Item(30,4).Value = "'0.09744"
$ws.Cells.Item(30,5).Value = "  -2.74%  "

# Row 31
$ws.Cells.Item(31,4).Value = "'1.488"
$ws.Cells.Item(31,5).Value = "  -1.09%  "

# Row 32
$ws.Cells.Item(32,4).Value = "'4.286"
$ws.Cells.Item(32,5).Value = "  +0.34%  "

# Row 33
$ws.Cells.Item(33,4).Value = "'4.172"
$ws.Cells.Item(33,5).Value = "  +1.08%  "

# Row 34
$ws.Cells.Item(34,4).Value = "'0.04885"
$ws.Cells.Item(34,5).Value = "  +1.18%  "

# Row 35
$ws.Cells.Item(35,4).Value = "'1.126"
$ws.Cells.Item(35,5).Value = "  +0.62%  "

# Row 36
$ws.Cells.Item(36,4).Value = "'0.6974"
$ws.Cells.Item(36,5).Value = "  +0.36%  "

# Row 37
$ws.Cells.Item(37,4).Value = "'2.734"
$ws.Cells.Item(37,5).Value = "  +0.80%  "

# Row 38
$ws.Cells.Item(38,4).Value = "'0.01894"
$ws.Cells.Item(38,5).Value = "  +2.22%  "

# Row 39
$ws.Cells.Item(39,4).Value = "'2.836"
$ws.Cells.Item(39,5).Value = "  +3.48%  "

# Row 40
$ws.Cells.Item(40,4).Value = "'76.09"
$ws.Cells.Item(40,5).Value = "  +4.52%  "

# Row 41
$ws.Cells.Item(41,4).Value = "'6.302"
$ws.Cells.Item(41,5).Value = "  -0.10%  "

# Row 42
$ws.Cells.Item(42,4).Value = "'2.003"
$ws.Cells.Item(42,5).Value = "  +1.37%  "

# Row 43
$ws.Cells.Item(43,4).Value = "'0.4264"
$ws.Cells.Item(43,5).Value = "  +1.72%  "

# Row 44
$ws.Cells.Item(44,5).Value = "  +0.26%  "

# Row 45
$ws.Cells.Item(45,4).Value = "'0.8366"
$ws.Cells.Item(45,5).Value = "  -0.08%  "

# Row 46
$ws.Cells.Item(46,4).Value = "'101.73"
$ws.Cells.Item(46,5).Value = "  -0.07%  "

# Row 47
$ws.Cells.Item(47,4).Value = "'9.613"
$ws.Cells.Item(47,5).Value = "  +3.03%  "

# Row 48
$ws.Cells.Item(48,2).Value = "Elrond"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(48,4).Value = "'35.50"
$ws.Cells.Item(48,5).Value = "  +0.09%  "

# Row 49
$ws.Cells.Item(49,2).Value = "Aptos"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(49,4).Value = "'6.998"
$ws.Cells.Item(49,5).Value = "  +0.27%  "

# Row 50
$ws.Cells.Item(50,4).Value = "'903.67"
$ws.Cells.Item(50,5).Value = "  -1.92%  "

# Row 51
$ws.Cells.Item(51,4).Value = "'0.05771"
$ws.Cells.Item(51,5).Value = "  +2.40%  "
